$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update author name (merged A2:B2)
$ws.Range("A2").Value = "Antoine Rochat"

# Row 5 - update activity text and hours
$ws.Range("B5").Value = "Présentation des directives du projet, constitution des groupes et discussion des différentes idées de projet. Au fin de la séance : décision finale de la proposition : BlaajjPaint"
$ws.Range("C5").Value = 1.5
$ws.Rows("5:5").RowHeight = 44.65

# Row 6 - change date, activity text, hours
$ws.Range("A6").Value = 43153
$ws.Range("B6").Value = "Relecture et complétion de la proposition à rendre"
$ws.Range("C6").Value = 1
$ws.Rows("6:6").RowHeight = 18.4

# Row 7 - new entry
$ws.Range("A7").Value = 43157
$ws.Range("B7").Value = "Retour du professeur sur notre proposition de projet et discussion entre l'équipe concernant les dates de rencontres et les fonctionnalités à mettre dans le cahier des charges"
$ws.Range("C7").Value = 1.5
$ws.Rows("7:7").RowHeight = 42.75

# Row 8 - new entry
$ws.Range("A8").Value = 43158
$ws.Range("B8").Value = "Réunion du groupe afin de discuter des souhaits de chacun, des spécificités du projets et des fonctionnalités (générales ou optionnelles). Conception également d'un schéma de dépendances fonctionnelles afin de faciliter le futur diagramme de Gantt"
$ws.Range("C8").Value = 4.5
$ws.Rows("8:8").RowHeight = 71.25

# Row 9 - new entry
$ws.Range("A9").Value = 43161
$ws.Range("B9").Value = "Rédaction et discussion sur le cahier des charges ainsi que début de conception du diagramme de Gantt et répartition des heures"
$ws.Range("C9").Value = 2
$ws.Rows("9:9").RowHeight = 28.5

# Update selection to match new state (active cell A2, selection A2:B2)
$ws.Range("A2:B2").Select() | Out-Null

$wb.Save()
